$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Re-apply the text number format to the header row and the "Row" label
# column (A), mirroring the style bump seen for these cells.
$ws.Range("A1:C1").NumberFormat = "@"
$ws.Range("A2:A14").NumberFormat = "@"

# Update the predicted-distance values in column B with the refreshed
# figures from the latest run (ful-path.csv).
$ws.Range("B2").Value  = 28724.416059566604
$ws.Range("B3").Value  = 96498.875125477236
$ws.Range("B4").Value  = 96484.654994095006
$ws.Range("B5").Value  = 81949.788452993496
$ws.Range("B6").Value  = 48389.937403419724
$ws.Range("B7").Value  = 104549.86787078655
$ws.Range("B8").Value  = 102228.83713563389
$ws.Range("B9").Value  = 81280.101526164945
$ws.Range("B10").Value = 82870.820477338653
$ws.Range("B11").Value = 87521.314502545225
$ws.Range("B12").Value = 87480.365440534282
$ws.Range("B13").Value = 121875.69746683838
$ws.Range("B14").Value = 50698.799158333859
